$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Lab/author labels (column F) ---
# Written in this order (F1, F3, F2) so the shared-string table is built in
# the same order as the authored workbook: "Lab 4, Part 2", "esv and cjd41",
# "Edward Venator and Chris Dickey".
$ws.Range("F1").Value = "Lab 4, Part 2"
$ws.Range("F3").Value = "esv and cjd41"
$ws.Range("F2").Value = "Edward Venator and Chris Dickey"

# --- Print orientation: landscape ---
$ws.PageSetup.Orientation = 2

# --- Move/resize the collector-current chart down below the data table ---
$chartObj = $ws.ChartObjects(1)
$chartObj.Left = 6.75
$chartObj.Top = 122.62496062992126
$chartObj.Width = 746.3096490526575
$chartObj.Height = 370.8751181102362

# --- Selection moves to F3 (matches the author's last-saved cursor spot) ---
$ws.Range("F3").Select()
